# Subnetting of all port.xlsx -- "Edited Subneting of WAN's"
# The /30 WAN subnet mask notes were wrong (used the Host-bit value,
# 255.255.255.254, as the subnet mask). Correct them to the real /30
# subnet mask, 255.255.255.252, for WAN-1, WAN-2 and WAN-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

$ws.Range("E17:E19").Value = "/30`n255.255.255.252"

# Leave the selection where the author ended up after the edit.
$ws.Range("E23").Select() | Out-Null
